# Update the "Valores" values report with recalculated counts
# (fix for competencia/ano function when generating reports)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values matrix, rows 1-16, columns A-L
$data = @(
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(11,0,0,7,0,7,0,5,0,12,6,9),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(4,0,0,1,0,1,0,1,0,4,2,1),
    @(7,0,0,4,0,3,0,5,0,2,2,2),
    @(1,0,0,0,0,0,0,0,0,2,1,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(1,0,0,0,0,1,0,0,0,2,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(7,0,0,4,0,5,0,0,0,8,5,8),
    @(31,0,0,16,0,17,0,11,0,30,16,20)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}
